$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 183; this shifts existing rows 183:227 down to 184:228
# and keeps the number format (s="2") on column D inherited from the row above.
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row 183 with the new weekly record
# (same Mercado/Region/Categoria/Variedad/Calidad/Unidad/Origen/Clasificacion
# as every other row in this sheet).
$ws.Range("A183").Value = 10
$ws.Range("B183").Value = "Vega Modelo de Temuco"
$ws.Range("C183").Value = "La Araucanía"
$ws.Range("D183").Value = 44785
$ws.Range("E183").Value = 9
$ws.Range("F183").Value = 100112005
$ws.Range("G183").Value = "Puerro"
$ws.Range("H183").Value = "Azul de Maquehue"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 30
$ws.Range("K183").Value = 15000
$ws.Range("L183").Value = 15000
$ws.Range("M183").Value = 15000
$ws.Range("N183").Value = "$/docena de paquetes"
$ws.Range("O183").Value = "Provincia de Cautín"
$ws.Range("P183").Value = 1250
$ws.Range("Q183").Value = 12
$ws.Range("R183").Value = "Hortaliza"
